$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project 1")

# Row 6: new "Grid Traveller" entry (filled in the same order the
# original author typed them, so new shared-string slots line up)
$ws.Range("B6").Value = "DP, Recursion"
$ws.Range("C6").Value = "Grid Traveller"
$ws.Range("E6").Value = "O(n*m)"
$ws.Range("F6").Value = "Identify the base cases. First solve the problem recursively(brute force) then add hash table to store the values."
$ws.Range("D6").Value = "https://leetcode.com/problems/unique-paths/"

# Row 5: add the missing Link cell (D5) for the Fibonacci row
$ws.Range("D5").Value = "https://leetcode.com/problems/fibonacci-number/"

# Match the styling used for the other Link column cell (D4): small 8pt
# Arial font with wrapped text
$ws.Range("D5:D6").Font.Size = 8
$ws.Range("D5:D6").WrapText = $true

# Update the active selection to D5, matching the authored change
[void]$ws.Range("D5").Select()
